# Update cryptos list - Fri Mar 8 22:25:36 UTC 2024
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    # Force the cell to Text format before assigning so that Excel doesn't
    # auto-convert numeric-looking strings (e.g. "144.80" -> 144.8) and
    # the written value keeps its literal/verbatim text content.
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $val
}

# Row 2 - Bitcoin
$ws.Range("D2").Value = "68.335.22"
$ws.Range("E2").Value = "  +1.49%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "3.906.57"
$ws.Range("E3").Value = "  +0.70%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.07%  "

# Row 5 - BNB
Set-TextValue "D5" "481.21"
$ws.Range("E5").Value = "  +2.10%  "

# Row 6 - Solana
Set-TextValue "D6" "144.80"
$ws.Range("E6").Value = "  -0.62%  "

# Row 7 - XRP
Set-TextValue "D7" "0.621"
$ws.Range("E7").Value = "  -2.24%  "

# Row 8 - USDC
Set-TextValue "D8" "0.998"
$ws.Range("E8").Value = "  -0.08%  "

# Row 9 - Cardano
Set-TextValue "D9" "0.722"
$ws.Range("E9").Value = "  -3.44%  "

# Row 10 - Dogecoin
$ws.Range("E10").Value = "  +6.81%  "

# Row 11 - ShibaInu
Set-TextValue "D11" "0.0000351"
$ws.Range("E11").Value = "  +12.51%  "

# Row 12 - Avalanche
Set-TextValue "D12" "42.53"
$ws.Range("E12").Value = "  -2.31%  "

# Row 13 - Polkadot
Set-TextValue "D13" "10.51"
$ws.Range("E13").Value = "  +0.33%  "

# Row 14 - WrappedliquidstakedEther2.0
$ws.Range("D14").Value = "4.530.31"
$ws.Range("E14").Value = "  +0.51%  "

# Row 15 - was WrappedEther, now Uniswap
$ws.Range("B15").Value = "Uniswap"
$ws.Range("C15").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
Set-TextValue "D15" "14.59"
$ws.Range("E15").Value = "  -1.76%  "

# Row 16 - was Uniswap, now WrappedEther
$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D16").Value = "3.923.65"
$ws.Range("E16").Value = "  +1.64%  "

# Row 17 - TRON
$ws.Range("E17").Value = "  -0.33%  "

# Row 18 - Chainlink
Set-TextValue "D18" "19.67"
$ws.Range("E18").Value = "  -2.25%  "

# Row 19 - Polygon
Set-TextValue "D19" "1.13"
$ws.Range("E19").Value = "  -3.38%  "

# Row 20 - WrappedBTC
$ws.Range("D20").Value = "68.422.51"
$ws.Range("E20").Value = "  +1.18%  "

# Row 21 - BitcoinCash
Set-TextValue "D21" "434.55"
$ws.Range("E21").Value = "  -0.65%  "

# Row 22 - was InternetComputer(DFINITY), now ImmutableX
$ws.Range("B22").Value = "ImmutableX"
$ws.Range("C22").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextValue "D22" "3.37"
$ws.Range("E22").Value = "  +1.91%  "

# Row 23 - was ImmutableX, now InternetComputer(DFINITY)
$ws.Range("B23").Value = "InternetComputer(DFINITY)"
$ws.Range("C23").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextValue "D23" "14.61"
$ws.Range("E23").Value = "  -2.20%  "

# Row 24 - Litecoin
Set-TextValue "D24" "87.43"
$ws.Range("E24").Value = "  -2.25%  "

# Row 25 - Filecoin
Set-TextValue "D25" "11.70"
$ws.Range("E25").Value = "  +17.39%  "

# Row 26 - PancakeSwap
Set-TextValue "D26" "3.58"
$ws.Range("E26").Value = "  -1.17%  "

# Row 27 - EthereumClassic
Set-TextValue "D27" "38.04"
$ws.Range("E27").Value = "  -0.24%  "

# Row 28 - RenderToken
Set-TextValue "D28" "10.40"
$ws.Range("E28").Value = "  +3.18%  "

# Row 29 - LEO
Set-TextValue "D29" "5.84"
$ws.Range("E29").Value = "  +6.35%  "

# Row 30 - Bittensor
Set-TextValue "D30" "700.88"
$ws.Range("E30").Value = "  -4.32%  "

# Row 31 - was Hedera, now Cosmos
$ws.Range("B31").Value = "Cosmos"
$ws.Range("C31").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-TextValue "D31" "13.32"
$ws.Range("E31").Value = "  -4.02%  "

# Row 32 - was Cosmos, now Hedera
$ws.Range("B32").Value = "Hedera"
$ws.Range("C32").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue "D32" "0.130"
$ws.Range("E32").Value = "  -3.05%  "

# Row 33 - Toncoin
Set-TextValue "D33" "2.87"
$ws.Range("E33").Value = "  +3.31%  "

# Row 34 - PEPE
$ws.Range("D34").Value = "0.0₃0921"
$ws.Range("E34").Value = "  +33.90%  "

# Row 35 - InjectiveProtocol
Set-TextValue "D35" "41.42"
$ws.Range("E35").Value = "  -6.79%  "

# Row 36 - OKB
Set-TextValue "D36" "59.37"
$ws.Range("E36").Value = "  +2.20%  "

# Row 37 - NEARProtocol
Set-TextValue "D37" "5.69"
$ws.Range("E37").Value = "  +2.81%  "

# Row 38 - Kaspa
$ws.Range("E38").Value = "  -7.63%  "

# Row 39 - Dai
Set-TextValue "D39" "0.998"
$ws.Range("E39").Value = "  -0.27%  "

# Row 40 - VeChain
Set-TextValue "D40" "0.0473"
$ws.Range("E40").Value = "  -2.59%  "

# Row 41 - WEMIXToken
$ws.Range("E41").Value = "  +10.78%  "

# Row 42 - Fetch.AI
Set-TextValue "D42" "2.75"
$ws.Range("E42").Value = "  +7.68%  "

# Row 43 - ThetaToken
Set-TextValue "D43" "2.99"
$ws.Range("E43").Value = "  +2.31%  "

# Row 44 - was Stellar, now TheGraph
$ws.Range("B44").Value = "TheGraph"
$ws.Range("C44").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
Set-TextValue "D44" "0.340"
$ws.Range("E44").Value = "  -2.59%  "

# Row 45 - was TheGraph, now Stellar
$ws.Range("B45").Value = "Stellar"
$ws.Range("C45").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextValue "D45" "0.141"
$ws.Range("E45").Value = "  -0.46%  "

# Row 46 - FirstDigitalUSD
$ws.Range("E46").Value = "  -0.15%  "

# Row 47 - LidoDAOToken
Set-TextValue "D47" "3.42"
$ws.Range("E47").Value = "  -1.46%  "

# Row 48 - ARBITRUM
Set-TextValue "D48" "2.14"
$ws.Range("E48").Value = "  -1.34%  "

# Row 49 - Monero
Set-TextValue "D49" "146.03"
$ws.Range("E49").Value = "  +1.20%  "

# Row 50 - ApeXProtocol
$ws.Range("E50").Value = "  -4.37%  "

# Row 51 - Stacks
Set-TextValue "D51" "2.84"
$ws.Range("E51").Value = "  -1.97%  "
